# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - linked from the slide master  ("Integral" palette)
#   ppt/theme/theme2.xml  - linked from the notes master   ("Office Theme" palette)
# The target edit swaps their contents, so that the slide master's theme
# becomes the "Office Theme" palette (and the notes master's theme becomes
# the "Integral" palette). The font scheme and format scheme are identical
# between the two themes, so the only real content change is the 12-colour
# colour scheme (and the theme's display name).
#
# The PowerPoint object model exposes the *slide master's* colour scheme as
# SlideMaster.ColorScheme.Colors(i).RGB (i = 1..12, in clrScheme document
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Recolour it with
# the "Office Theme" palette that currently lives in theme2.xml.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function ToComRgb([int]$r, [int]$g, [int]$b) {
    # PowerPoint's RGB colour values are packed as 0x00BBGGRR.
    return $r + ($g * 256) + ($b * 65536)
}

# "Office Theme" colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeTheme = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $rgb = $officeTheme[$i - 1]
    $colorScheme.Colors($i).RGB = ToComRgb $rgb[0] $rgb[1] $rgb[2]
}
